$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("E1").Value = "ideal_age"
$ws.Range("F1").Value = "ideal_gender"

# ideal_age values (E2:E31) and ideal_gender values (F2:F31)
$idealAge = @(18,21,30,60,32,40,19,27,33,31,20,18,45,38,16,20,45,47,37,31,26,19,21,17,40,26,49,33,27,27)
$idealGender = @("male","female","female","female","female","male","male","male","female","female","male","male","male","male","female","male","female","female","female","male","female","female","male","male","male","female","female","female","female","male")

for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $idealAge[$i]
    $ws.Cells.Item($row, 6).Value = $idealGender[$i]
}

$ws.Range("E31").Select()
